# Generate Report for Handoff
#
# The CI run produced a fresh handoff for e2e\b.md: it is now "Ready for
# handoff" with a new xliff + timestamp, and the previously out-of-date
# handback file on "a.md" now carries an explanatory error message that got
# copied onto b's row as well (mirrors the upstream localization-status
# report generator). Column P (Error Detail) also needed to widen to fit the
# new long text.

$wb = $excel.ActiveWorkbook

$newHandoffFile  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5009cfa63e0387e635b520ae2e5d597218f74476/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/794a843001e9dfe85d645740b1612b698794b4aa/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is e2e\b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 16:38:24"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text (otherwise "False" is auto-coerced to a
# Boolean, like real Excel); re-applying the Normal style afterwards drops
# the quote-prefix formatting flag that the apostrophe trick sets.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 16:38:20"
$wsZhCn.Range("P3").Value = $newHandoffFile
# ColumnWidth is in "characters"; Excel stores the column width in the file
# with ~5px of cell padding added, so asking for an even 40 on-disk means
# requesting a slightly smaller character width here.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 16:38:24"
$wsDeDe.Range("P3").Value = $newHandoffFile
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
